$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "60.936.67"; E = "  -1.11%  " },
    @{ Row = 3; D = "3.374.82"; E = "  -0.18%  " },
    @{ Row = 4; E = "  -0.08%  " },
    @{ Row = 5; D = "571.97"; E = "  -0.96%  " },
    @{ Row = 6; D = "136.28"; E = "  -0.28%  " },
    @{ Row = 7; E = "  +0.06%  " },
    @{ Row = 8; D = "3.373.27"; E = "  -0.15%  " },
    @{ Row = 9; D = "0.468"; E = "  -1.28%  " },
    @{ Row = 10; E = "  +2.12%  " },
    @{ Row = 11; D = "0.121"; E = "  -3.24%  " },
    @{ Row = 12; D = "0.379"; E = "  -2.68%  " },
    @{ Row = 13; D = "3.951.50"; E = "  -0.24%  " },
    @{ Row = 14; E = "  +0.51%  " },
    @{ Row = 15; D = "25.64"; E = "  +0.54%  " },
    @{ Row = 16; B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "3.373.54"; E = "  -0.38%  " },
    @{ Row = 17; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D = "0.0000170"; E = "  -3.35%  " },
    @{ Row = 18; D = "61.121.84"; E = "  -1.04%  " },
    @{ Row = 19; D = "13.77"; E = "  -2.75%  " },
    @{ Row = 20; D = "5.73"; E = "  -1.24%  " },
    @{ Row = 21; D = "9.29"; E = "  -2.02%  " },
    @{ Row = 22; D = "375.20"; E = "  -1.04%  " },
    @{ Row = 23; D = "3.517.59"; E = "  -0.30%  " },
    @{ Row = 24; D = "0.546"; E = "  -2.52%  " },
    @{ Row = 25; D = "0.999"; E = "  +0.16%  " },
    @{ Row = 26; D = "71.01"; E = "  -0.21%  " },
    @{ Row = 27; D = "0.0000123"; E = "  -0.87%  " },
    @{ Row = 28; D = "1.64"; E = "  -5.54%  " },
    @{ Row = 29; D = "0.176"; E = "  +10.50%  " },
    @{ Row = 30; D = "0.999" },
    @{ Row = 31; D = "7.37"; E = "  -2.97%  " },
    @{ Row = 32; D = "8.03"; E = "  -1.81%  " },
    @{ Row = 33; E = "  -1.60%  " },
    @{ Row = 34; E = "  -0.05%  " },
    @{ Row = 35; D = "23.34"; E = "  -0.15%  " },
    @{ Row = 36; D = "5.10"; E = "  -4.56%  " },
    @{ Row = 37; D = "1.54"; E = "  -1.00%  " },
    @{ Row = 38; D = "6.77"; E = "  -0.95%  " },
    @{ Row = 39; D = "164.65"; E = "  +0.21%  " },
    @{ Row = 40; D = "0.0758"; E = "  -3.49%  " },
    @{ Row = 41; E = "  -0.08%  " },
    @{ Row = 42; B = "Mantle"; C = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; D = "0.775"; E = "  -0.58%  " },
    @{ Row = 43; B = "EnergySwap"; C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D = "25.09"; E = "  +1.06%  " },
    @{ Row = 44; D = "1.68"; E = "  -2.53%  " },
    @{ Row = 45; D = "4.31"; E = "  -2.29%  " },
    @{ Row = 46; D = "1.17"; E = "  -4.47%  " },
    @{ Row = 47; D = "2.535.04"; E = "  +8.41%  " },
    @{ Row = 48; D = "6.75"; E = "  -1.45%  " },
    @{ Row = 49; D = "22.83"; E = "  +0.09%  " },
    @{ Row = 50; D = "2.42"; E = "  +3.53%  " },
    @{ Row = 51; D = "0.0258"; E = "  -1.48%  " }
)

foreach ($item in $updates) {
    $r = $item.Row
    if ($item.ContainsKey("B")) { $ws.Range("B$r").Value = $item.B }
    if ($item.ContainsKey("C")) { $ws.Range("C$r").Value = $item.C }
    if ($item.ContainsKey("D")) {
        $cellD = $ws.Range("D$r")
        $cellD.NumberFormat = "@"
        $cellD.Value = $item.D
    }
    if ($item.ContainsKey("E")) { $ws.Range("E$r").Value = $item.E }
}

Write-Output "Applied $($updates.Count) row updates"